$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (preserve inline/shared-string type,
# avoid Excel auto-converting numeric-looking strings like "420.48" into
# a Number cell). Sets text number-format, writes the value, then resets
# the cell style back to Normal so no stray style id is left behind.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Update Price (D) and Volume(1h) (E) columns ---
Set-TextValue "D2" "66.775.90"
$ws.Range("E2").Value = "  +2.66%  "

Set-TextValue "D3" "3.720.75"
$ws.Range("E3").Value = "  +6.22%  "

$ws.Range("E4").Value = "  +0.17%  "

Set-TextValue "D5" "420.48"
$ws.Range("E5").Value = "  +0.21%  "

Set-TextValue "D6" "131.97"
$ws.Range("E6").Value = "  +0.18%  "

Set-TextValue "D7" "3.711.38"
$ws.Range("E7").Value = "  +6.18%  "

Set-TextValue "D8" "0.645"
$ws.Range("E8").Value = "  -0.74%  "

$ws.Range("E9").Value = "  +0.04%  "

Set-TextValue "D10" "0.771"
$ws.Range("E10").Value = "  -0.57%  "

Set-TextValue "D11" "0.183"
$ws.Range("E11").Value = "  +13.57%  "

Set-TextValue "D12" "0.0000404"
$ws.Range("E12").Value = "  +54.85%  "

Set-TextValue "D13" "42.91"
$ws.Range("E13").Value = "  -0.36%  "

Set-TextValue "D14" "10.55"
$ws.Range("E14").Value = "  +7.24%  "

Set-TextValue "D15" "4.298.51"
$ws.Range("E15").Value = "  +5.87%  "

$ws.Range("E16").Value = "  -0.83%  "

Set-TextValue "D17" "20.79"
$ws.Range("E17").Value = "  +1.96%  "

Set-TextValue "D18" "3.702.97"
$ws.Range("E18").Value = "  +5.85%  "

Set-TextValue "D19" "13.23"
$ws.Range("E19").Value = "  +6.00%  "

$ws.Range("E20").Value = "  +4.20%  "

Set-TextValue "D21" "66.838.61"
$ws.Range("E21").Value = "  +2.85%  "

Set-TextValue "D22" "446.55"
$ws.Range("E22").Value = "  -2.93%  "

Set-TextValue "D23" "16.52"
$ws.Range("E23").Value = "  +24.54%  "

Set-TextValue "D24" "89.75"
$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("E25").Value = "  -0.92%  "

Set-TextValue "D26" "38.11"
$ws.Range("E26").Value = "  +12.30%  "

Set-TextValue "D27" "10.22"
$ws.Range("E27").Value = "  +2.68%  "

Set-TextValue "D28" "3.33"
$ws.Range("E28").Value = "  -0.59%  "

Set-TextValue "D29" "5.07"
$ws.Range("E29").Value = "  +4.19%  "

Set-TextValue "D30" "12.82"
$ws.Range("E30").Value = "  +2.74%  "

Set-TextValue "D31" "0.124"
$ws.Range("E31").Value = "  +9.76%  "

$ws.Range("E32").Value = "  +2.39%  "

$ws.Range("E33").Value = "  -2.99%  "

Set-TextValue "D34" "0.165"
$ws.Range("E34").Value = "  +1.10%  "

Set-TextValue "D35" "41.84"
$ws.Range("E35").Value = "  +5.05%  "

Set-TextValue "D36" "57.20"
$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("E37").Value = "  +0.02%  "

Set-TextValue "D38" "0.0496"
$ws.Range("E38").Value = "  -1.23%  "

Set-TextValue "D39" "0.0₃0747"
$ws.Range("E39").Value = "  +5.25%  "

Set-TextValue "D40" "3.07"
$ws.Range("E40").Value = "  +31.84%  "

$ws.Range("E41").Value = "  +0.96%  "

Set-TextValue "D42" "28.60"
$ws.Range("E42").Value = "  +31.47%  "

Set-TextValue "D43" "0.998"
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("E44").Value = "  +4.89%  "

Set-TextValue "D45" "3.25"
$ws.Range("E45").Value = "  +33.20%  "

Set-TextValue "D46" "147.64"
$ws.Range("E46").Value = "  +0.82%  "

$ws.Range("E47").Value = "  +5.84%  "

Set-TextValue "D50" "2.90"
$ws.Range("E50").Value = "  -5.79%  "

Set-TextValue "D51" "0.309"
$ws.Range("E51").Value = "  -1.45%  "

# --- Rows 48/49 swap: NEARProtocol and WEMIXToken swap rank positions ---
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D48" "2.67"
$ws.Range("E48").Value = "  -3.98%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D49" "4.38"
$ws.Range("E49").Value = "  -2.68%  "

